$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1518.0454
$ws.Range("I40").Value = 1352.5
$ws.Range("J40").Value = 1656
$ws.Range("K40").Value = 1352.5
$ws.Range("L40").Value = 1656
$ws.Range("M40").Value = -1177.5
$ws.Range("N40").Value = -2006
$ws.Range("H125").Value = 3572.5881
$ws.Range("J125").Value = 3572.5881
$ws.Range("L125").Value = 32153.2929
$ws.Range("N125").Value = -37073.2929
$ws.Range("H137").Value = 2167.775
$ws.Range("I137").Value = 1588.3529
$ws.Range("J137").Value = 2596.0435
$ws.Range("K137").Value = 4765.0587
$ws.Range("L137").Value = 7788.130500000001
$ws.Range("M137").Value = -2215.0587
$ws.Range("N137").Value = -12888.1305

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 10006
$ws.Range("I16").Value = 10006
$ws.Range("K16").Value = 10006
$ws.Range("M16").Value = -9719
$ws.Range("H102").Value = 4091
$ws.Range("I102").Value = 4212.222
$ws.Range("K102").Value = 4212.222
$ws.Range("M102").Value = -2590.222
$ws.Range("H110").Value = 1906.875
$ws.Range("I110").Value = 2101.111
$ws.Range("J110").Value = 1657.1428
$ws.Range("K110").Value = 2101.111
$ws.Range("L110").Value = 1657.1428
$ws.Range("M110").Value = -56.11099999999988
$ws.Range("N110").Value = -5747.1428

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 919.8182
$ws.Range("I94").Value = 914.875
$ws.Range("J94").Value = 933
$ws.Range("K94").Value = 914.875
$ws.Range("L94").Value = 933
$ws.Range("M94").Value = -463.875
$ws.Range("N94").Value = -1835
$ws.Range("H105").Value = 2930.8235
$ws.Range("I105").Value = 2921.2727
$ws.Range("J105").Value = 2948.3333
$ws.Range("K105").Value = 2921.2727
$ws.Range("L105").Value = 2948.3333
$ws.Range("M105").Value = -1174.2727
$ws.Range("N105").Value = -6442.3333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1900
$ws.Range("I16").Value = 1900
$ws.Range("J16").Value = 1900
$ws.Range("K16").Value = 1900
$ws.Range("L16").Value = 1900
$ws.Range("M16").Value = -1613
$ws.Range("N16").Value = -2474
$ws.Range("H22").Value = 440.77777
$ws.Range("I22").Value = 263.14285
$ws.Range("J22").Value = 1062.5
$ws.Range("K22").Value = 263.14285
$ws.Range("L22").Value = 1062.5
$ws.Range("M22").Value = 86.85714999999999
$ws.Range("N22").Value = -1762.5
$ws.Range("H94").Value = 1162
$ws.Range("I94").Value = 794.8
$ws.Range("J94").Value = 1303.2307
$ws.Range("K94").Value = 794.8
$ws.Range("L94").Value = 1303.2307
$ws.Range("M94").Value = -343.8
$ws.Range("N94").Value = -2205.2307
$ws.Range("H105").Value = 2098.625
$ws.Range("I105").Value = 1756.6666
$ws.Range("J105").Value = 2303.8
$ws.Range("K105").Value = 1756.6666
$ws.Range("L105").Value = 2303.8
$ws.Range("M105").Value = -9.666600000000017
$ws.Range("N105").Value = -5797.8
$ws.Range("H113").Value = 1900
$ws.Range("I113").Value = 1900
$ws.Range("J113").Value = 1900
$ws.Range("K113").Value = 1900
$ws.Range("L113").Value = 1900
$ws.Range("M113").Value = 270
$ws.Range("N113").Value = -6240

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 2662.5
$ws.Range("I40").Value = 216.66667
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 866.66668
$ws.Range("L40").Value = 40000
$ws.Range("M40").Value = -797.66668
$ws.Range("N40").Value = -40138
$ws.Range("H55").Value = 6931.6665
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 6931.6665
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 20794.9995
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -21148.9995
$ws.Range("H68").Value = 1182
$ws.Range("I68").Value = 881.7143
$ws.Range("J68").Value = 1357.1666
$ws.Range("K68").Value = 2645.1429
$ws.Range("L68").Value = 4071.4998
$ws.Range("M68").Value = -1834.1429
$ws.Range("N68").Value = -5693.4998
$ws.Range("H71").Value = 1182
$ws.Range("I71").Value = 881.7143
$ws.Range("J71").Value = 1357.1666
$ws.Range("K71").Value = 7935.428699999999
$ws.Range("L71").Value = 12214.4994
$ws.Range("M71").Value = -3879.428699999999
$ws.Range("N71").Value = -20326.4994
$ws.Range("H106").Value = 7488.3335
$ws.Range("J106").Value = 7488.3335
$ws.Range("L106").Value = 22465.0005
$ws.Range("N106").Value = -24357.0005
$ws.Range("H109").Value = 2007.5
$ws.Range("I109").Value = 1000
$ws.Range("K109").Value = 3000
$ws.Range("M109").Value = -1960

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 21333.334
$ws.Range("H94").Value = 30448
$ws.Range("J94").Value = 30448
$ws.Range("L94").Value = 30448
$ws.Range("N94").Value = -31800
$ws.Range("H97").Value = 1313.1154
$ws.Range("I97").Value = 1262.8572
$ws.Range("J97").Value = 1524.2
$ws.Range("K97").Value = 1262.8572
$ws.Range("L97").Value = 1524.2
$ws.Range("M97").Value = -766.8571999999999
$ws.Range("N97").Value = -2516.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1459.6
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H27").Value = 1459.6
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H61").Value = 15381.8125
$ws.Range("J61").Value = 3233.3333
$ws.Range("L61").Value = 3233.3333
$ws.Range("N61").Value = -3637.3333
$ws.Range("H93").Value = 2990
$ws.Range("I93").Value = 2990
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2990
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -1742
$ws.Range("N93").ClearContents()
$ws.Range("H113").Value = 15381.8125
$ws.Range("J113").Value = 3233.3333
$ws.Range("L113").Value = 3233.3333
$ws.Range("N113").Value = -7573.3333

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I47").Value = 50000
$ws.Range("K47").Value = 50000
$ws.Range("M47").Value = -49428
$ws.Range("H96").Value = 1500
$ws.Range("I96").Value = 1500
$ws.Range("K96").Value = 1500
$ws.Range("M96").Value = -127
$ws.Range("H113").Value = 1084.5294
$ws.Range("I113").Value = 712.8333
$ws.Range("J113").Value = 1976.6
$ws.Range("K113").Value = 2138.4999
$ws.Range("L113").Value = 5929.799999999999
$ws.Range("M113").Value = 31.5001000000002
$ws.Range("N113").Value = -10269.8
